# Word COM-interop script applying the edit described by the diff.
#
# The document's real (visible-text) change is in the heading that used
# to read:
#     "A backendel kommunikáló service réteg elkészítése"
# and now reads:
#     "A backenddel kommunikáló service réteg a desktop alkalmazásban"
#
# (All the other hunks in the diff are cosmetic: XML-namespace / Word
# "mc:Ignorable" list churn, w:proofErr (spell/grammar-check) marker
# shuffling and run-splitting that leaves the underlying text identical,
# and the internal "_GoBack" bookmark that Word silently re-positions to
# mark the place of the last edit - none of these change what a reader
# actually sees, so there is nothing for Find/Replace to do for them.)

$d = $word.ActiveDocument

# 1) "backendel" -> "backenddel" (typo fix: missing second "d").
$d.Content.Find.Execute("backendel", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "backenddel", 2)

# 2) "elkészítése" -> "a desktop alkalmazásban" (re-worded ending of the
#    heading, clarifying that this is about the service layer inside the
#    desktop application).
$d.Content.Find.Execute("elkészítése", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a desktop alkalmazásban", 2)
